# Update the dSF (column F) values to reflect the repulled / recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F4").Value = 1
$ws.Range("F8").Value = 2
$ws.Range("F11").Value = -5
$ws.Range("F14").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 0
